$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# Capture existing "Title" header/value before we shift columns
$title = $ws.Range("D1").Value2
$titleVal = $ws.Range("D2").Value2

# Insert new axis-arrow / percentage columns (D,E,F), pushing the old
# Title column out to G. Assignment order matters for how new strings
# are appended to the shared string table, so keep this exact order.
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"
$ws.Range("A2").Value = "Ol"
$ws.Range("E2").Value = "Orthopyroxene (%)"
$ws.Range("D2").Value = "Olivine (%)"
$ws.Range("F2").Value = "Clinopyroxene (%)"
$ws.Range("G1").Value = $title
$ws.Range("G2").Value = $titleVal

# Widen the new columns to match their content
$ws.Range("D1:F1").ColumnWidth = 13.25

# Make "axes" the active sheet, with F8 selected
$ws.Activate()
$ws.Range("F8").Select()
